# Append the new daily APR data point (row 9) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "2025-09-20T11:39"
$ws.Range("C9").Value = 1.6866366606922056
